# Updated warning messages for attorney.
# Applies the changes to the "Case_Data" sheet (Sheet1):
#  - Clears the stray empty cell at G27
#  - Appends four new case rows (29-32)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the empty inline-string cell that used to live at G27
$ws.Range("G27").ClearContents()

# New row 29
$ws.Cells.Item(29, 1).Value = "21TRD09386"
$ws.Cells.Item(29, 2).Value = "Bunner"
$ws.Cells.Item(29, 3).Value = "DUS UCM"
$ws.Cells.Item(29, 4).Value = "'4510.111"
$ws.Cells.Item(29, 5).Value = "UCM"
$ws.Cells.Item(29, 6).Value = "No Contest"
$ws.Cells.Item(29, 7).Value = "Guilty"
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = "'0"

# New row 30
$ws.Cells.Item(30, 1).Value = "21TRD09386"
$ws.Cells.Item(30, 2).Value = "Bunner"
$ws.Cells.Item(30, 3).Value = "TAIL LIGHTS-REAR LICENSE PLATE"
$ws.Cells.Item(30, 4).Value = "'4513.05"
$ws.Cells.Item(30, 5).Value = "MM"
$ws.Cells.Item(30, 6).Value = "No Contest"
$ws.Cells.Item(30, 7).Value = "Guilty"
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = "'0"

# New row 31
$ws.Cells.Item(31, 1).Value = "21CRB01268"
$ws.Cells.Item(31, 2).Value = "Bunner"
$ws.Cells.Item(31, 3).Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(31, 4).Value = "2925.14(C)"
$ws.Cells.Item(31, 5).Value = "M4"
$ws.Cells.Item(31, 6).Value = "No Contest"
$ws.Cells.Item(31, 7).Value = "Guilty"
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = "'0"

# New row 32
$ws.Cells.Item(32, 1).Value = "21CRB01268"
$ws.Cells.Item(32, 2).Value = "Hemmeter"
$ws.Cells.Item(32, 3).Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(32, 4).Value = "2925.14(C)"
$ws.Cells.Item(32, 5).Value = "M4"
$ws.Cells.Item(32, 6).Value = "No Contest"
$ws.Cells.Item(32, 7).Value = "Guilty"
$ws.Cells.Item(32, 8).Value = "'50"
$ws.Cells.Item(32, 9).Value = "'25"
